$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Excel_vs_ML")
$ws3 = $wb.Worksheets.Item("Exec_Summary")

# Row 2
$ws1.Range("E2").Value = 46060
$ws1.Range("F2").Value = 46059

# Row 3
$ws1.Range("E3").Value = 46060
$ws1.Range("F3").Value = 46059
$ws1.Range("H3").Value = 449512.12
$ws1.Range("I3").Value = 44
$ws1.Range("J3").Value = 12
$ws1.Range("K3").Value = 434181.49
$ws1.Range("L3").Value = 103.53
$ws1.Range("M3").Value = 103082.51
$ws1.Range("N3").Value = 8590.209999999999
$ws1.Range("O3").Value = "On Track"
$ws1.Range("P3").Value = "Overdelivered"
$ws1.Range("Q3").Value = 103082.51
$ws1.Range("R3").Value = "YES"

# Row 4
$ws1.Range("E4").Value = 46060
$ws1.Range("F4").Value = 46059

# Row 5
$ws1.Range("E5").Value = 46060
$ws1.Range("F5").Value = 46059

# Row 6
$ws1.Range("E6").Value = 46060
$ws1.Range("F6").Value = 46059

# Row 7
$ws1.Range("E7").Value = 46060
$ws1.Range("F7").Value = 46059

# Row 8
$ws1.Range("E8").Value = 46060
$ws1.Range("F8").Value = 46059
$ws1.Range("H8").Value = 34460.47
$ws1.Range("I8").Value = 16
$ws1.Range("J8").Value = 43
$ws1.Range("K8").Value = 37109.4
$ws1.Range("L8").Value = 92.86
$ws1.Range("M8").Value = 102380.45
$ws1.Range("N8").Value = 2380.94
$ws1.Range("O8").Value = "Underpacing"

# Row 9
$ws1.Range("E9").Value = 46060
$ws1.Range("F9").Value = 46059

# Row 10
$ws1.Range("E10").Value = 46060
$ws1.Range("F10").Value = 46059
$ws1.Range("H10").Value = 90191.59
$ws1.Range("I10").Value = 13
$ws1.Range("J10").Value = 44
$ws1.Range("K10").Value = 102371.16
$ws1.Range("L10").Value = 88.09999999999999
$ws1.Range("M10").Value = 358666.56
$ws1.Range("N10").Value = 8151.51
$ws1.Range("O10").Value = "Underpacing"
$ws1.Range("R10").Value = "NO"

# Row 11
$ws1.Range("E11").Value = 46060
$ws1.Range("F11").Value = 46059

# Row 12
$ws1.Range("E12").Value = 46060
$ws1.Range("F12").Value = 46059
$ws1.Range("H12").Value = 209956.21
$ws1.Range("I12").Value = 36
$ws1.Range("J12").Value = 46
$ws1.Range("K12").Value = 245768.81
$ws1.Range("L12").Value = 85.43000000000001
$ws1.Range("M12").Value = 349850.52
$ws1.Range("N12").Value = 7605.45

# Row 13
$ws1.Range("E13").Value = 46060
$ws1.Range("F13").Value = 46059
$ws1.Range("H13").Value = 48203.28
$ws1.Range("I13").Value = 22
$ws1.Range("J13").Value = 59
$ws1.Range("K13").Value = 46003.78
$ws1.Range("L13").Value = 104.78
$ws1.Range("M13").Value = 121174.26
$ws1.Range("N13").Value = 2053.8
$ws1.Range("P13").Value = "Overdelivered"
$ws1.Range("Q13").Value = 121174.26
$ws1.Range("R13").Value = "YES"

# Row 14
$ws1.Range("E14").Value = 46060
$ws1.Range("F14").Value = 46059

# Row 15
$ws1.Range("E15").Value = 46060
$ws1.Range("F15").Value = 46059

# Row 16
$ws1.Range("E16").Value = 46060
$ws1.Range("F16").Value = 46059

# Row 17
$ws1.Range("E17").Value = 46060
$ws1.Range("F17").Value = 46059

# Row 18
$ws1.Range("E18").Value = 46060
$ws1.Range("F18").Value = 46059

# Row 19
$ws1.Range("E19").Value = 46060
$ws1.Range("F19").Value = 46059

# Row 20
$ws1.Range("E20").Value = 46060
$ws1.Range("F20").Value = 46059
$ws1.Range("H20").Value = 156172.13
$ws1.Range("I20").Value = 37
$ws1.Range("J20").Value = 53
$ws1.Range("K20").Value = 174179.52
$ws1.Range("L20").Value = 89.66
$ws1.Range("M20").Value = 267507.78
$ws1.Range("N20").Value = 5047.32

# Row 21
$ws1.Range("E21").Value = 46060
$ws1.Range("F21").Value = 46059

# Row 22
$ws1.Range("E22").Value = 46060
$ws1.Range("F22").Value = 46059

# Row 23
$ws1.Range("E23").Value = 46060
$ws1.Range("F23").Value = 46059

# Row 24
$ws1.Range("E24").Value = 46060
$ws1.Range("F24").Value = 46059

# Row 25
$ws1.Range("E25").Value = 46060
$ws1.Range("F25").Value = 46059

# Row 26
$ws1.Range("E26").Value = 46060
$ws1.Range("F26").Value = 46059

# Row 27
$ws1.Range("E27").Value = 46060
$ws1.Range("F27").Value = 46059
$ws1.Range("H27").Value = 361627.27
$ws1.Range("I27").Value = 50
$ws1.Range("J27").Value = 33
$ws1.Range("K27").Value = 348077.64
$ws1.Range("L27").Value = 103.89
$ws1.Range("M27").Value = 216181.61
$ws1.Range("N27").Value = 6550.96
$ws1.Range("O27").Value = "On Track"
$ws1.Range("P27").Value = "Overdelivered"
$ws1.Range("Q27").Value = 216181.61
$ws1.Range("R27").Value = "YES"

# Row 28
$ws1.Range("E28").Value = 46060
$ws1.Range("F28").Value = 46059

# Row 29
$ws1.Range("E29").Value = 46060
$ws1.Range("F29").Value = 46059

# Row 30
$ws1.Range("E30").Value = 46060
$ws1.Range("F30").Value = 46059

# Row 31
$ws1.Range("E31").Value = 46060
$ws1.Range("F31").Value = 46059

# Row 32
$ws1.Range("E32").Value = 46060
$ws1.Range("F32").Value = 46059
$ws1.Range("H32").Value = 243681.7
$ws1.Range("I32").Value = 62
$ws1.Range("J32").Value = 58
$ws1.Range("K32").Value = 296949.14
$ws1.Range("L32").Value = 82.06
$ws1.Range("M32").Value = 331058.57
$ws1.Range("N32").Value = 5707.91

# Row 33
$ws1.Range("E33").Value = 46060
$ws1.Range("F33").Value = 46059
$ws1.Range("H33").Value = 54101.04
$ws1.Range("I33").Value = 32
$ws1.Range("J33").Value = 34
$ws1.Range("K33").Value = 62718.53
$ws1.Range("L33").Value = 86.26000000000001
$ws1.Range("M33").Value = 75255.92999999999
$ws1.Range("N33").Value = 2213.41

# Row 34
$ws1.Range("E34").Value = 46060
$ws1.Range("F34").Value = 46059

# Row 35
$ws1.Range("E35").Value = 46060
$ws1.Range("F35").Value = 46059

# Row 36
$ws1.Range("E36").Value = 46060
$ws1.Range("F36").Value = 46059

# Row 37
$ws1.Range("E37").Value = 46060
$ws1.Range("F37").Value = 46059

# Row 38
$ws1.Range("E38").Value = 46060
$ws1.Range("F38").Value = 46059

# Row 39
$ws1.Range("E39").Value = 46060
$ws1.Range("F39").Value = 46059

# Row 40
$ws1.Range("E40").Value = 46060
$ws1.Range("F40").Value = 46059
$ws1.Range("H40").Value = 46798.99
$ws1.Range("I40").Value = 37
$ws1.Range("J40").Value = 28
$ws1.Range("K40").Value = 68358.61
$ws1.Range("L40").Value = 68.45999999999999
$ws1.Range("M40").Value = 73290.46000000001
$ws1.Range("N40").Value = 2617.52

# Row 41
$ws1.Range("E41").Value = 46060
$ws1.Range("F41").Value = 46059

# Row 42
$ws1.Range("E42").Value = 46060
$ws1.Range("F42").Value = 46059

# Row 43
$ws1.Range("E43").Value = 46060
$ws1.Range("F43").Value = 46059

# Row 44
$ws1.Range("E44").Value = 46060
$ws1.Range("F44").Value = 46059

# Row 45
$ws1.Range("E45").Value = 46060
$ws1.Range("F45").Value = 46059
$ws1.Range("H45").Value = 251835.04
$ws1.Range("I45").Value = 56
$ws1.Range("J45").Value = 23
$ws1.Range("K45").Value = 251369.83
$ws1.Range("L45").Value = 100.19
$ws1.Range("M45").Value = 102775.97
$ws1.Range("N45").Value = 4468.52
$ws1.Range("O45").Value = "On Track"
$ws1.Range("P45").Value = "Overdelivered"
$ws1.Range("Q45").Value = 102775.97
$ws1.Range("R45").Value = "YES"

# Row 46
$ws1.Range("E46").Value = 46060
$ws1.Range("F46").Value = 46059

# Exec_Summary updates
$ws3.Range("B3").Value = 25
$ws3.Range("B4").Value = 10
$ws3.Range("B5").Value = 4
$ws3.Range("B6").Value = 543214.35
